## feat: add 2022-Q4 data
##
## 1. Insert a new worksheet named "2022-Q4" right before the existing
##    "2022-Q3" sheet (so the tab order becomes 总计, 2022-Q4, 2022-Q3,
##    2021-Q2, 2021-Q1, 2020-Q4).
## 2. Populate "2022-Q4" with its fund-holding table (13 rows).
## 3. Insert a new summary row for "2022-Q4" at the top of the "总计"
##    sheet's data (row 2), pushing the other quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: insert + name the new sheet
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)          # currently "2022-Q3"
$newSheet = $wb.Worksheets.Add($beforeSheet)   # inserted before it
$newSheet.Name = "2022-Q4"

# After the insert, sheet order/index is:
#   1 总计, 2 2022-Q4 (new), 3 2022-Q3, 4 2021-Q2, 5 2021-Q1, 6 2020-Q4
# Re-resolve template sheets fresh (by index) rather than reusing the
# $beforeSheet handle, since indices shifted.
$headerTemplate = $wb.Worksheets.Item(3)   # "2022-Q3" - exact header match
$indexTemplate  = $wb.Worksheets.Item(5)   # "2021-Q1" - has >= 13 data rows

# Copy the header row (style + text both already correct for every
# column, including D1 = "基金规模").
$headerTemplate.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

# Copy the bold/bordered index-column style for the 13 data rows.
$indexTemplate.Range("A2:A14").Copy($newSheet.Range("A2:A14"))

function Set-TextCell($cell, $val) {
    # Force text storage (otherwise Excel auto-coerces numeric-looking
    # strings like "4.65" into real numbers), then drop back to the
    # default "Normal" style so no stray number-format survives.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-FundRow($ws, $row, $aIdx, $code, $name, $scale, $pos, $ratio, $mv, $rank) {
    $ws.Cells.Item($row, 1).Value = $aIdx
    Set-TextCell $ws.Cells.Item($row, 2) $code
    Set-TextCell $ws.Cells.Item($row, 3) $name
    Set-TextCell $ws.Cells.Item($row, 4) $scale
    Set-TextCell $ws.Cells.Item($row, 5) $pos
    Set-TextCell $ws.Cells.Item($row, 6) $ratio
    Set-TextCell $ws.Cells.Item($row, 7) $mv
    $ws.Cells.Item($row, 8).Value = $rank
}

# ---------------------------------------------------------------------
# Step 2: fill the "2022-Q4" fund table
# ---------------------------------------------------------------------
Set-FundRow $newSheet 2  0 "550001" "信诚四季红混合"             "4.65" "82.79" "2.67" "0.1242" 10
Set-FundRow $newSheet 3  1 "519013" "海富通风格优势混合"         "3.23" "92.73" "2.86" "0.0924" 4
Set-FundRow $newSheet 4  2 "013051" "汇泉臻心致远混合A"          "2.02" "76.67" "2.67" "0.0539" 7
Set-FundRow $newSheet 5  3 "013052" "汇泉臻心致远混合C"          "1.37" "76.67" "2.67" "0.0366" 7
Set-FundRow $newSheet 6  4 "006973" "太平睿盈混合A"              "3.17" "29.39" "1.05" "0.0333" 10
Set-FundRow $newSheet 7  5 "011284" "中信保诚龙腾精选混合"       "1.08" "83.70" "2.67" "0.0288" 10
Set-FundRow $newSheet 8  6 "007669" "太平睿盈混合C"              "0.99" "29.39" "1.05" "0.0104" 10
Set-FundRow $newSheet 9  7 "015201" "创金合信动态平衡混合C"      "0.21" "67.50" "4.19" "0.0088" 1
Set-FundRow $newSheet 10 8 "015200" "创金合信动态平衡混合A"      "0.19" "67.50" "4.19" "0.0080" 1
Set-FundRow $newSheet 11 9 "011438" "红塔红土盛昌优选混合A"      "0.13" "92.67" "4.56" "0.0059" 8
Set-FundRow $newSheet 12 10 "001412" "德邦鑫星价值灵活配置混合A" "0.13" "35.79" "1.83" "0.0024" 10
Set-FundRow $newSheet 13 11 "011439" "红塔红土盛昌优选混合C"     "0.04" "92.67" "4.56" "0.0018" 8
Set-FundRow $newSheet 14 12 "002112" "德邦鑫星价值灵活配置混合C" "0.02" "35.79" "1.83" "0.0004" 10

# ---------------------------------------------------------------------
# Step 3: insert the "2022-Q4" summary row into "总计" (sheet 1)
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

$totalSheet.Rows.Item(2).Insert()

# Clean up the bare formatting the row-insert carried over from row 1
# (the header) so the new data cells end up with no explicit style,
# matching the rest of the table.
$totalSheet.Range("B2:D2").ClearFormats()

# A2 needs the same bold/bordered index-column style as A3:A6; grab it
# from A3 (which still carries the original style after the insert).
$totalSheet.Cells.Item(3, 1).Copy($totalSheet.Cells.Item(2, 1))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 13
$totalSheet.Cells.Item(2, 4).Value = 0.41

# Renumber the shifted rows' index column (0,1,2,3 -> 1,2,3,4) so the
# whole A column reads 0,1,2,3,4 top to bottom.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
